# Add more solved problems to the "All Questions" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Questions")

# New LeetCode problems solved (appended right after the existing table, rows 17-19)
$rows = @(
    @{ Num = 16; Topic = "Two Pointers"; Url = "https://leetcode.com/problems/merge-sorted-array/" },
    @{ Num = 17; Topic = "Binary Tree";  Url = "https://leetcode.com/problems/binary-tree-inorder-traversal/" },
    @{ Num = 18; Topic = "Binary Tree";  Url = "https://leetcode.com/problems/same-tree/" }
)

$r = 17
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row.Num
    $ws.Range("B$r").Value = $row.Topic
    $ws.Range("C$r").Value = $row.Url
    $ws.Range("D$r").Value = "Easy"
    $ws.Range("E$r").Value = "NA"
    $ws.Hyperlinks.Add($ws.Range("C$r"), $row.Url, "", "", $row.Url)

    # Hyperlinks.Add re-formats the cell with the built-in blue/underlined
    # "Hyperlink" look; restore the plain row formatting used by the rest
    # of the table by copying formats from the row directly above.
    $ws.Range("A16:E16").Copy()
    $ws.Range("A$r" + ":E$r").PasteSpecial(-4122)

    $r = $r + 1
}
$excel.CutCopyMode = $false

# Update the saved view state: scroll "All Questions" down and select E22
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("E22").Select()

# Reset the lingering multi-area selections on the other sheets to a plain A1 selection
foreach ($name in @("Strategy", "Topic data", "GFG")) {
    $other = $wb.Worksheets.Item($name)
    $other.Range("A1").Select()
}

$ws.Select()
